# Fruta / hortaliza, semanal
# Inserts a new weekly data row for "Femacal de La Calera" (Arándano (blue))
# right after row 67, shifting all subsequent rows down by one
# (old row 68 -> new row 69, ..., old row 99 -> new row 100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68 (pushes rows 68..99 down to 69..100)
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record
$ws.Cells.Item(68, 1).Value = 3
$ws.Cells.Item(68, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(68, 3).Value = "Coquimbo"
$ws.Cells.Item(68, 4).Value = 44466
$ws.Cells.Item(68, 5).Value = 5
$ws.Cells.Item(68, 6).Value = "Fruta"
$ws.Cells.Item(68, 7).Value = 100101
$ws.Cells.Item(68, 8).Value = "Berries"
$ws.Cells.Item(68, 9).Value = 100101001
$ws.Cells.Item(68, 10).Value = "Arándano (blue)"
$ws.Cells.Item(68, 11).Value = "Sin especificar"
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 45
$ws.Cells.Item(68, 14).Value = 10000
$ws.Cells.Item(68, 15).Value = 10000
$ws.Cells.Item(68, 16).Value = 10000
$ws.Cells.Item(68, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(68, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(68, 19).Value = 6667
$ws.Cells.Item(68, 20).Value = 1.5
